# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with freshly scraped figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "30.641.53"; Volume = "+1.30%" },
    @{ Row = 3; Price = "1.860.54"; Volume = "+0.16%" },
    @{ Row = 4; Price = $null; Volume = "+0.16%" },
    @{ Row = 5; Price = "235.03"; Volume = "+0.94%" },
    @{ Row = 6; Price = "1.001"; Volume = "+0.10%" },
    @{ Row = 7; Price = "0.4709"; Volume = "-0.25%" },
    @{ Row = 8; Price = "0.2760"; Volume = "+0.96%" },
    @{ Row = 9; Price = "0.06357"; Volume = "-1.02%" },
    @{ Row = 10; Price = "17.55"; Volume = "+8.23%" },
    @{ Row = 11; Price = "1.844.72"; Volume = "+0.31%" },
    @{ Row = 12; Price = "0.07443"; Volume = "-0.07%" },
    @{ Row = 13; Price = "5.178"; Volume = "+3.64%" },
    @{ Row = 14; Price = "84.92"; Volume = "-0.53%" },
    @{ Row = 15; Price = "0.6318"; Volume = "+0.29%" },
    @{ Row = 16; Price = "30.625.97"; Volume = "+1.43%" },
    @{ Row = 17; Price = "242.96"; Volume = "+4.37%" },
    @{ Row = 18; Price = "1.001"; Volume = "+0.10%" },
    @{ Row = 19; Price = "12.81"; Volume = "+0.71%" },
    @{ Row = 20; Price = "0.000007367"; Volume = "+0.14%" },
    @{ Row = 21; Price = $null; Volume = "+0.26%" },
    @{ Row = 22; Price = "4.995"; Volume = "-0.76%" },
    @{ Row = 23; Price = "6.050"; Volume = "+0.85%" },
    @{ Row = 24; Price = "9.345"; Volume = "+1.08%" },
    @{ Row = 25; Price = "165.30"; Volume = "+0.22%" },
    @{ Row = 26; Price = "18.10"; Volume = "+1.28%" },
    @{ Row = 27; Price = "1.888"; Volume = "+0.13%" },
    @{ Row = 28; Price = "0.1017"; Volume = "-0.66%" },
    @{ Row = 29; Price = "1.383"; Volume = "+0.16%" },
    @{ Row = 30; Price = "4.067"; Volume = "-1.44%" },
    @{ Row = 31; Price = "3.865"; Volume = "-1.60%" },
    @{ Row = 32; Price = "0.04930"; Volume = "+0.80%" },
    @{ Row = 33; Price = $null; Volume = "+0.49%" },
    @{ Row = 34; Price = $null; Volume = "-2.17%" },
    @{ Row = 35; Price = "2.713"; Volume = "+0.71%" },
    @{ Row = 36; Price = "0.01912"; Volume = "+0.75%" },
    @{ Row = 37; Price = "2.686"; Volume = "+1.88%" },
    @{ Row = 38; Price = "0.8808"; Volume = "-2.51%" },
    @{ Row = 39; Price = "1.995"; Volume = "+1.32%" },
    @{ Row = 40; Price = "105.34"; Volume = "-0.35%" },
    @{ Row = 41; Price = "1.001"; Volume = "+0.21%" },
    @{ Row = 42; Price = "5.546"; Volume = "+0.50%" },
    @{ Row = 43; Price = "0.4078"; Volume = "-0.62%" },
    @{ Row = 44; Price = "7.275"; Volume = "+2.55%" },
    @{ Row = 45; Price = "63.41"; Volume = "+3.86%" },
    @{ Row = 46; Price = "0.1215"; Volume = "+1.60%" },
    @{ Row = 47; Price = "33.59"; Volume = "+0.95%" },
    @{ Row = 48; Price = "8.578"; Volume = "-1.67%" },
    @{ Row = 49; Price = "0.05549"; Volume = "-0.33%" },
    @{ Row = 50; Price = "1.373"; Volume = "-2.35%" },
    @{ Row = 51; Price = "0.3694"; Volume = "-0.01%" }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D$row")
        if ($u.Price -match '^-?[0-9]+(\.[0-9]+)?$') {
            # A plain decimal would be auto-coerced to a number by Excel and
            # lose formatting (e.g. trailing zeros); force text with a quote
            # prefix, then strip the resulting style back to Normal.
            $priceCell.Value = "'" + $u.Price
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $u.Price
        }
    }
    $ws.Range("E$row").Value = "  " + $u.Volume + "  "
}
